$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.311.56"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.089.80"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.14"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.40"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.081.75"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("E10").Value = "  +5.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.60"
$ws.Range("E11").Value = "  -2.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  -1.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.06"
$ws.Range("E14").Value = "  +3.85%  "

$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.601.72"
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.195.15"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.096.32"
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.03"
$ws.Range("E19").Value = "  -2.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.15"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.37"
$ws.Range("E23").Value = "  -2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.97"
$ws.Range("E24").Value = "  -1.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.87"
$ws.Range("E25").Value = "  -3.29%  "

$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.00"
$ws.Range("E28").Value = "  +8.29%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.65"
$ws.Range("E30").Value = "  -1.11%  "

$ws.Range("E31").Value = "  -2.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.88"
$ws.Range("E32").Value = "  +0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.52"
$ws.Range("E34").Value = "  -1.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0841"
$ws.Range("E35").Value = "  -3.44%  "

$ws.Range("E36").Value = "  +1.62%  "

$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("E38").Value = "  -5.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.96"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.14"
$ws.Range("E40").Value = "  -1.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "430.94"
$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.68"
$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.884.00"
$ws.Range("E43").Value = "  -1.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0364"
$ws.Range("E44").Value = "  -1.56%  "

$ws.Range("E45").Value = "  -3.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.107"
$ws.Range("E46").Value = "  -3.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.78"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.43"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.86"
$ws.Range("E51").Value = "  -2.92%  "
